$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.09
$ws.Range("G2").Value = 970
$ws.Range("H2").Value = 1.01
$ws.Range("I2").Value = 970
$ws.Range("J2").Value = 1.09
$ws.Range("V2").Value = 1.09

# Row 3
$ws.Range("I3").Value = 2.4
$ws.Range("K3").Value = 3.65
$ws.Range("M3").Value = 1.07
$ws.Range("S3").Value = 3.35
$ws.Range("AI3").Value = 65
$ws.Range("AN3").Value = 95
$ws.Range("AO3").Value = 34

# Row 4
$ws.Range("U4").Value = 1.92

# Row 5
$ws.Range("F5").Value = 1.09
$ws.Range("G5").Value = 600
$ws.Range("J5").Value = 1.09
$ws.Range("V5").Value = 1.13
$ws.Range("W5").Value = 1.17

# Row 6
$ws.Range("F6").Value = 2.96
$ws.Range("G6").Value = 3.6
$ws.Range("J6").Value = 2.96
$ws.Range("M6").Value = 1.08
$ws.Range("W6").Value = 1.39

# Row 7
$ws.Range("J7").Value = 3.2
$ws.Range("K7").Value = 3.95
$ws.Range("T7").Value = 1.76
$ws.Range("V7").Value = 1.33

# Row 8
$ws.Range("G8").Value = 7.6
$ws.Range("L8").Value = 1.23
$ws.Range("N8").Value = 2.28
$ws.Range("S8").Value = 2.52
$ws.Range("T8").Value = 1.68
$ws.Range("U8").Value = 1.04
$ws.Range("W8").Value = 1.15

# Row 9
$ws.Range("F9").Value = 2.18
$ws.Range("G9").Value = 2.48
$ws.Range("H9").Value = 2.84
$ws.Range("J9").Value = 3.8
$ws.Range("L9").Value = 1.01
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 3.45
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 2.3
$ws.Range("Q9").Value = 1.59
$ws.Range("R9").Value = 1.52
$ws.Range("S9").Value = 2.48
$ws.Range("T9").Value = 1.56
$ws.Range("U9").Value = 2.42
$ws.Range("V9").Value = 1.46
$ws.Range("W9").Value = 1.67
$ws.Range("X9").Value = 24
$ws.Range("Y9").Value = 960
$ws.Range("Z9").Value = 30
$ws.Range("AA9").Value = 60
$ws.Range("AB9").Value = 960
$ws.Range("AC9").Value = 960
$ws.Range("AD9").Value = 960
$ws.Range("AE9").Value = 38
$ws.Range("AF9").Value = 960
$ws.Range("AG9").Value = 960
$ws.Range("AH9").Value = 960
$ws.Range("AI9").Value = 42
$ws.Range("AJ9").Value = 32
$ws.Range("AK9").Value = 23
$ws.Range("AL9").Value = 32
$ws.Range("AM9").Value = 75
$ws.Range("AN9").Value = 960
$ws.Range("AO9").Value = 26

# Row 10
$ws.Range("H10").Value = 1.69
$ws.Range("R10").Value = 1.52
$ws.Range("S10").Value = 2.82
$ws.Range("T10").Value = 1.74
$ws.Range("U10").Value = 2.24
$ws.Range("X10").Value = 20
$ws.Range("Y10").Value = 10.5
$ws.Range("Z10").Value = 11.5
$ws.Range("AA10").Value = 17.5
$ws.Range("AB10").Value = 22
$ws.Range("AD10").Value = 9.6
$ws.Range("AE10").Value = 16
$ws.Range("AF10").Value = 44
$ws.Range("AG10").Value = 20
$ws.Range("AH10").Value = 18
$ws.Range("AI10").Value = 28
$ws.Range("AJ10").Value = 130
$ws.Range("AK10").Value = 65
$ws.Range("AN10").Value = 55

# Row 11
$ws.Range("J11").Value = 3.5
$ws.Range("K11").Value = 3.55
$ws.Range("L11").Value = 1.43
$ws.Range("O11").Value = 1.34
$ws.Range("P11").Value = 1.94
$ws.Range("S11").Value = 3.6
$ws.Range("T11").Value = 1.8
$ws.Range("U11").Value = 2.14
$ws.Range("W11").Value = 1.81
$ws.Range("X11").Value = 13
$ws.Range("Y11").Value = 14
$ws.Range("Z11").Value = 25
$ws.Range("AA11").Value = 70
$ws.Range("AB11").Value = 9.6
$ws.Range("AC11").Value = 7.6
$ws.Range("AD11").Value = 15
$ws.Range("AE11").Value = 44
$ws.Range("AF11").Value = 13.5
$ws.Range("AG11").Value = 10.5
$ws.Range("AH11").Value = 17
$ws.Range("AI11").Value = 55
$ws.Range("AJ11").Value = 28
$ws.Range("AK11").Value = 23
$ws.Range("AL11").Value = 38
$ws.Range("AN11").Value = 17.5
$ws.Range("AO11").Value = 44

# Row 12
$ws.Range("H12").Value = 4.1
$ws.Range("I12").Value = 4.3
$ws.Range("J12").Value = 4.5
$ws.Range("K12").Value = 4.6
$ws.Range("P12").Value = 3.1
$ws.Range("Q12").Value = 1.44
$ws.Range("T12").Value = 1.47
$ws.Range("U12").Value = 2.84
$ws.Range("X12").Value = 32
$ws.Range("Z12").Value = 38
$ws.Range("AA12").Value = 85
$ws.Range("AB12").Value = 17.5
$ws.Range("AD12").Value = 18
$ws.Range("AE12").Value = 40
$ws.Range("AG12").Value = 11
$ws.Range("AH12").Value = 15
$ws.Range("AI12").Value = 38
$ws.Range("AK12").Value = 15.5
$ws.Range("AL12").Value = 23
$ws.Range("AM12").Value = 50
$ws.Range("AO12").Value = 23

# Row 13
$ws.Range("F13").Value = 1.77
$ws.Range("G13").Value = 1.78
$ws.Range("J13").Value = 4.2
$ws.Range("K13").Value = 4.3
$ws.Range("L13").Value = 1.28
$ws.Range("Q13").Value = 1.78
$ws.Range("R13").Value = 1.48
$ws.Range("S13").Value = 2.98
$ws.Range("T13").Value = 1.77
$ws.Range("U13").Value = 2.2
$ws.Range("W13").Value = 2.28
$ws.Range("X13").Value = 18.5
$ws.Range("Y13").Value = 20
$ws.Range("Z13").Value = 40
$ws.Range("AA13").Value = 120
$ws.Range("AC13").Value = 9.199999999999999
$ws.Range("AD13").Value = 19.5
$ws.Range("AE13").Value = 60
$ws.Range("AF13").Value = 10.5
$ws.Range("AG13").Value = 9.6
$ws.Range("AH13").Value = 18
$ws.Range("AI13").Value = 60
$ws.Range("AL13").Value = 29
$ws.Range("AM13").Value = 90
$ws.Range("AO13").Value = 55

# Row 14
$ws.Range("H14").Value = 4.1
$ws.Range("J14").Value = 3.05
$ws.Range("Q14").Value = 2.08
$ws.Range("S14").Value = 3.75
$ws.Range("V14").Value = 1.23

# Row 16
$ws.Range("Q16").Value = 1.87
